# Fix ADS8686 default values, all must be hex values of the form 0xVALUE
#
# The "Default Value" column (C) on the ADS8686 sheet stored a handful of
# register defaults as plain numbers (0 or 2) instead of the hex-string
# convention ("0x00" / "0x02") used everywhere else in that column. Rewrite
# those cells so every Default Value on the sheet is a proper 0xVALUE string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADS8686")

# Rows whose Default Value (column C) was the bare number 0 -> "0x00"
$rowsToZeroHex = @(2,3,11,13,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44)
foreach ($r in $rowsToZeroHex) {
    $ws.Range("C$r").Value = "0x00"
}

# Row whose Default Value (column C) was the bare number 2 -> "0x02"
$ws.Range("C12").Value = "0x02"
